# Update column G ("K") values for rows 2-17 with newly computed strike counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 2
    6  = 2
    7  = 1
    8  = 2
    9  = 1
    10 = 1
    11 = 1
    12 = 2
    13 = 2
    14 = 2
    15 = 1
    16 = 0
    17 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
